# Add 23 new "pelada" score rows (198-220) to the "Jogadores" sheet,
# mirroring the author's upload of additional Saturday-game results.
# Column layout (row 1 headers): A=Jogadores, B=Pontos, C=Vitorias,
# D=Empate, E=Derrotas, F=Gols, G=Partidas, H=Tarde de Vitoria,
# I=La barca, J=Craque do Dia, K=Gols Sofridos.  Column B is left blank,
# matching every pre-existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{Row=198; Name="Leandrinho";  Vals=@(4,3,2,1,1,0,0,0,0)},
  @{Row=199; Name="Corinthiano"; Vals=@(4,3,2,1,1,0,0,0,0)},
  @{Row=200; Name="Juscielio";   Vals=@(4,3,2,2,1,0,0,0,0)},
  @{Row=201; Name="Bruno";       Vals=@(4,3,2,0,1,0,0,0,0)},
  @{Row=202; Name="Cabeleira";   Vals=@(4,3,2,4,1,0,0,0,0)},
  @{Row=203; Name="David";       Vals=@(3,0,5,2,1,0,0,0,0)},
  @{Row=204; Name="Boneco";      Vals=@(3,0,5,1,1,0,0,0,0)},
  @{Row=205; Name="Romario";     Vals=@(3,0,5,2,1,0,0,0,0)},
  @{Row=206; Name="Marlin";      Vals=@(3,0,5,1,1,0,0,0,0)},
  @{Row=207; Name="Ismael";      Vals=@(3,0,5,0,1,0,0,0,0)},
  @{Row=208; Name="Athos";       Vals=@(5,2,3,7,1,1,0,1,0)},
  @{Row=209; Name="Ranyeri";     Vals=@(5,2,3,1,1,1,0,0,0)},
  @{Row=210; Name="Marcelão";    Vals=@(5,2,3,3,1,1,0,0,0)},
  @{Row=211; Name="Digão";       Vals=@(5,2,3,1,1,1,0,0,0)},
  @{Row=212; Name="Jorge";       Vals=@(5,2,3,0,1,1,0,0,0)},
  @{Row=213; Name="Eder";        Vals=@(2,1,4,0,1,0,1,0,0)},
  @{Row=214; Name="Adriano";     Vals=@(2,1,4,1,1,0,1,0,0)},
  @{Row=215; Name="Joãozinho";   Vals=@(2,1,4,2,1,0,1,0,0)},
  @{Row=216; Name="Marcos";      Vals=@(2,1,4,1,1,0,1,0,0)},
  @{Row=217; Name="Du";          Vals=@(2,1,4,1,1,0,1,0,0)},
  @{Row=218; Name="Matheus";     Vals=@(7,3,2,0,1,1,0,0,6)},
  @{Row=219; Name="Alan";        Vals=@(5,1,5,0,1,0,0,0,8)},
  @{Row=220; Name="Chelin";      Vals=@(2,0,7,1,1,0,1,0,13)}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $col = 3
    foreach ($v in $r.Vals) {
        $ws.Cells.Item($r.Row, $col).Value = $v
        $col++
    }
}

# Move the active selection to match the author's final cursor position.
[void]$ws.Range("I220").Select()

Write-Output "Added rows 198-220"
